$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a brand-new worksheet "2022-Q4" right after "总计" (i.e. before
#    the worksheet that is currently named "2022-Q3"). All the other
#    quarterly tabs keep their own name + data; they simply shift one
#    position to the right in the tab order.
# ---------------------------------------------------------------------------
$oldQ3 = $wb.Worksheets.Item(2)
$q4 = $wb.Worksheets.Add($oldQ3, $null)
$q4.Name = "2022-Q4"

# Header row
$q4.Cells.Item(1,2).Value = "基金代码"
$q4.Cells.Item(1,3).Value = "基金名称"
$q4.Cells.Item(1,4).Value = "基金规模"
$q4.Cells.Item(1,5).Value = "股票总仓位"
$q4.Cells.Item(1,6).Value = "仓位占比"
$q4.Cells.Item(1,7).Value = "持有市值(亿元)"
$q4.Cells.Item(1,8).Value = "仓位排名"

# Row 2
$q4.Cells.Item(2,1).Value = 0
$q4.Cells.Item(2,2).Value = "012368"
$q4.Cells.Item(2,3).Value = "摩根士丹利华鑫优享臻选六个月持有期混合A"
$q4.Cells.Item(2,4).Value = "4.61"
$q4.Cells.Item(2,5).Value = "93.52"
$q4.Cells.Item(2,6).Value = "5.51"
$q4.Cells.Item(2,7).Value = "0.2540"
$q4.Cells.Item(2,8).Value = 8

# Row 3
$q4.Cells.Item(3,1).Value = 1
$q4.Cells.Item(3,2).Value = "233006"
$q4.Cells.Item(3,3).Value = "大摩领先优势混合"
$q4.Cells.Item(3,4).Value = "3.76"
$q4.Cells.Item(3,5).Value = "93.54"
$q4.Cells.Item(3,6).Value = "5.01"
$q4.Cells.Item(3,7).Value = "0.1884"
$q4.Cells.Item(3,8).Value = 8

# Row 4
$q4.Cells.Item(4,1).Value = 2
$q4.Cells.Item(4,2).Value = "000309"
$q4.Cells.Item(4,3).Value = "大摩品质生活精选股票"
$q4.Cells.Item(4,4).Value = "3.34"
$q4.Cells.Item(4,5).Value = "93.92"
$q4.Cells.Item(4,6).Value = "5.11"
$q4.Cells.Item(4,7).Value = "0.1707"
$q4.Cells.Item(4,8).Value = 8

# Row 5
$q4.Cells.Item(5,1).Value = 3
$q4.Cells.Item(5,2).Value = "010322"
$q4.Cells.Item(5,3).Value = "大摩新兴产业股票"
$q4.Cells.Item(5,4).Value = "2.06"
$q4.Cells.Item(5,5).Value = "93.45"
$q4.Cells.Item(5,6).Value = "5.27"
$q4.Cells.Item(5,7).Value = "0.1086"
$q4.Cells.Item(5,8).Value = 9

# Row 6
$q4.Cells.Item(6,1).Value = 4
$q4.Cells.Item(6,2).Value = "012369"
$q4.Cells.Item(6,3).Value = "摩根士丹利华鑫优享臻选六个月持有期混合C"
$q4.Cells.Item(6,4).Value = "0.31"
$q4.Cells.Item(6,5).Value = "93.52"
$q4.Cells.Item(6,6).Value = "5.51"
$q4.Cells.Item(6,7).Value = "0.0171"
$q4.Cells.Item(6,8).Value = 8

# ---------------------------------------------------------------------------
# 2. "总计" (summary) sheet: insert a new row 2 for the 2022-Q4 totals and
#    push the existing rows down by one (they keep their own values).
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$summary.Rows.Item(2).Insert()

$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q4"
$summary.Cells.Item(2,3).Value = 5
$summary.Cells.Item(2,4).Value = 0.74

# Renumber the index column (A) for the rows that shifted down.
$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(4,1).Value = 2
$summary.Cells.Item(5,1).Value = 3
$summary.Cells.Item(6,1).Value = 4
$summary.Cells.Item(7,1).Value = 5
